$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.01
$ws.Range("H3").Value = 4.5
$ws.Range("I3").Value = 6.6
$ws.Range("K3").Value = 6.4
$ws.Range("S3").Value = 3.95
$ws.Range("W3").Value = 1.96
$ws.Range("I5").Value = 1.53
$ws.Range("N5").Value = 4.5
$ws.Range("P5").Value = 2.22
$ws.Range("Q5").Value = 1.79
$ws.Range("T5").Value = 1.95
$ws.Range("V5").Value = 2.88
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 230
$ws.Range("AK5").Value = 110
$ws.Range("AN5").Value = 130
$ws.Range("H6").Value = 12.5
$ws.Range("I6").Value = 13.5
$ws.Range("J6").Value = 7.6
$ws.Range("K6").Value = 7.8
$ws.Range("N6").Value = 8
$ws.Range("P6").Value = 3.35
$ws.Range("Q6").Value = 1.39
$ws.Range("R6").Value = 1.93
$ws.Range("S6").Value = 2.02
$ws.Range("Y6").Value = 60
$ws.Range("AD6").Value = 46
$ws.Range("AE6").Value = 170
$ws.Range("H7").Value = 2.18
$ws.Range("I7").Value = 2.22
$ws.Range("N7").Value = 4.7
$ws.Range("Q7").Value = 1.74
$ws.Range("T7").Value = 1.64
$ws.Range("U7").Value = 2.44
$ws.Range("X7").Value = 18
$ws.Range("F8").Value = 1.45
$ws.Range("G8").Value = 1.47
$ws.Range("J8").Value = 5.2
$ws.Range("R8").Value = 1.44
$ws.Range("T8").Value = 2.08
$ws.Range("U8").Value = 1.87
$ws.Range("W8").Value = 3.1
$ws.Range("AA8").Value = 280
$ws.Range("AC8").Value = 11.5
$ws.Range("AJ8").Value = 12
$ws.Range("Q9").Value = 1.6
$ws.Range("AB9").Value = 19
$ws.Range("AN9").Value = 20
$ws.Range("F10").Value = 2.26
$ws.Range("G10").Value = 2.3
$ws.Range("H10").Value = 3.1
$ws.Range("R10").Value = 1.69
$ws.Range("Y10").Value = 21
$ws.Range("AC10").Value = 9.6
$ws.Range("F11").Value = 2.18
$ws.Range("G11").Value = 2.2
$ws.Range("H11").Value = 3.45
$ws.Range("I11").Value = 3.55
$ws.Range("K11").Value = 3.9
$ws.Range("S11").Value = 2.78
$ws.Range("V11").Value = 1.39
$ws.Range("W11").Value = 1.83
$ws.Range("Z11").Value = 26
$ws.Range("AC11").Value = 8.800000000000001
$ws.Range("AE11").Value = 36
$ws.Range("AF11").Value = 15
$ws.Range("AI11").Value = 40
$ws.Range("AO11").Value = 27
$ws.Range("I12").Value = 20
$ws.Range("J12").Value = 9.800000000000001
$ws.Range("P12").Value = 4.4
$ws.Range("S12").Value = 1.68
$ws.Range("T12").Value = 1.83
$ws.Range("AB12").Value = 18.5
$ws.Range("AD12").Value = 1000
$ws.Range("AN12").Value = 2.48
$ws.Range("F13").Value = 3.1
$ws.Range("H13").Value = 2.38
$ws.Range("I13").Value = 2.4
$ws.Range("N13").Value = 4.9
$ws.Range("R13").Value = 1.53
$ws.Range("S13").Value = 2.76
$ws.Range("V13").Value = 1.71
